$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1013.1818
$ws.Range("J97").Value = 1027.2222
$ws.Range("L97").Value = 3081.6666
$ws.Range("N97").Value = -4073.6666
$ws.Range("H138").Value = 4178.2324
$ws.Range("I138").Value = 3546.0952
$ws.Range("J138").Value = 4382.4614
$ws.Range("K138").Value = 10638.2856
$ws.Range("L138").Value = 13147.3842
$ws.Range("M138").Value = -5498.285600000001
$ws.Range("N138").Value = -23427.3842

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11004620
$ws.Range("I32").Value = 13711481
$ws.Range("J32").Value = 26794.445
$ws.Range("K32").Value = 13711481
$ws.Range("L32").Value = 26794.445
$ws.Range("M32").Value = -13711194
$ws.Range("N32").Value = -27368.445
$ws.Range("H61").Value = 13892924
$ws.Range("I61").Value = 27780500
$ws.Range("K61").Value = 27780500
$ws.Range("M61").Value = -27780288
$ws.Range("H132").Value = 2752363.8
$ws.Range("I132").Value = 7034.933
$ws.Range("J132").Value = 5920051
$ws.Range("K132").Value = 21104.799
$ws.Range("L132").Value = 17760153
$ws.Range("M132").Value = -18574.799
$ws.Range("N132").Value = -17765213
$ws.Range("H136").Value = 13892924
$ws.Range("I136").Value = 27780500
$ws.Range("K136").Value = 83341500
$ws.Range("M136").Value = -83338950

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31251744
$ws.Range("I20").Value = 1567.1364
$ws.Range("J20").Value = 100002136
$ws.Range("K20").Value = 1567.1364
$ws.Range("L20").Value = 100002136
$ws.Range("M20").Value = -1320.1364
$ws.Range("N20").Value = -100002630
$ws.Range("H68").Value = 42000
$ws.Range("J68").Value = 42000
$ws.Range("L68").Value = 42000
$ws.Range("N68").Value = -43622
$ws.Range("H71").Value = 42000
$ws.Range("J71").Value = 42000
$ws.Range("L71").Value = 126000
$ws.Range("N71").Value = -134112
$ws.Range("H135").Value = 61788.57
$ws.Range("J135").Value = 61788.57
$ws.Range("L135").Value = 61788.57
$ws.Range("N135").Value = -71928.57000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8540.261
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 8540.261
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 8540.261
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -9130.261
$ws.Range("H34").Value = 8540.261
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 8540.261
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 8540.261
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8944.261
$ws.Range("H86").Value = 4187.5
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 3916.6667
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 3916.6667
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -6162.6667
$ws.Range("H89").Value = 4187.5
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 3916.6667
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 19583.3335
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -30815.3335
$ws.Range("H99").Value = 2382.879
$ws.Range("I99").Value = 2112.3635
$ws.Range("K99").Value = 2112.3635
$ws.Range("M99").Value = -614.3634999999999
$ws.Range("H126").Value = 2382.879
$ws.Range("I126").Value = 2112.3635
$ws.Range("K126").Value = 6337.0905
$ws.Range("M126").Value = -3867.0905
$ws.Range("H132").Value = 95241310
$ws.Range("I132").Value = 250005000
$ws.Range("J132").Value = 33335834
$ws.Range("K132").Value = 750015000
$ws.Range("L132").Value = 100007502
$ws.Range("M132").Value = -750012470
$ws.Range("N132").Value = -100012562

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11343334
$ws.Range("J11").Value = 530002
$ws.Range("L11").Value = 530002
$ws.Range("N11").Value = -530280
$ws.Range("H18").Value = 8500
$ws.Range("J18").Value = 8500
$ws.Range("L18").Value = 8500
$ws.Range("M18").Value = -9086
$ws.Range("H107").Value = 434.76923
$ws.Range("I107").Value = 368.45456
$ws.Range("J107").Value = 799.5
$ws.Range("K107").Value = 368.45456
$ws.Range("L107").Value = 799.5
$ws.Range("M107").Value = 1551.54544
$ws.Range("N107").Value = -4639.5
$ws.Range("H132").Value = 37043972
$ws.Range("I132").Value = 66675990
$ws.Range("J132").Value = 3943.6667
$ws.Range("K132").Value = 200027970
$ws.Range("L132").Value = 11831.0001
$ws.Range("M132").Value = -200025440
$ws.Range("N132").Value = -16891.0001
$ws.Range("H135").Value = 48712.855
$ws.Range("J135").Value = 48712.855
$ws.Range("L135").Value = 48712.855
$ws.Range("N135").Value = -58852.855

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H25").Value = 836008.3
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 836008.3
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 836008.3
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -836468.3
$ws.Range("H132").Value = 2830.9636
$ws.Range("I132").Value = 2302.9375
$ws.Range("J132").Value = 3565.6086
$ws.Range("K132").Value = 6908.8125
$ws.Range("L132").Value = 10696.8258
$ws.Range("M132").Value = -4378.8125
$ws.Range("N132").Value = -15756.8258
$ws.Range("H134").Value = 23307.25
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 27743
$ws.Range("K134").Value = 10000
$ws.Range("L134").Value = 27743
$ws.Range("M134").Value = -4930
$ws.Range("N134").Value = -37883
$ws.Range("H136").Value = 3705378.2
$ws.Range("I136").Value = 1172.1034
$ws.Range("J136").Value = 10419252
$ws.Range("K136").Value = 3516.3102
$ws.Range("L136").Value = 31257756
$ws.Range("M136").Value = -966.3101999999999
$ws.Range("N136").Value = -31262856

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 40789.445
$ws.Range("J70").Value = 40789.445
$ws.Range("L70").Value = 40789.445
$ws.Range("N70").Value = -41419.445
$ws.Range("H73").Value = 40789.445
$ws.Range("J73").Value = 40789.445
$ws.Range("L73").Value = 40789.445
$ws.Range("N73").Value = -42973.445
$ws.Range("H136").Value = 3543.6223
$ws.Range("I136").Value = 4092.6428
$ws.Range("J136").Value = 2639.353
$ws.Range("K136").Value = 12277.9284
$ws.Range("L136").Value = 7918.059
$ws.Range("M136").Value = -9727.928400000001
$ws.Range("N136").Value = -13018.059
